# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# sheets to reflect refreshed counts from the data source.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 1849
    $ws.Range("F4").Value = 149
    $ws.Range("F6").Value = 6304
    $ws.Range("F7").Value = 168
}
